$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("S1").Value = "zugelassen"
$ws.Range("U1").Value = "umr-fb12-finalthesis_1.0.1_2025-02-06"
$ws.Range("Q1").Value = (Get-Date -Year 2000 -Month 12 -Day 1).Date

$ws.Range("A1:K1").Select()
